$d = $word.ActiveDocument

$replacements = @(
    @{old = "74×87="; new = "29×63="},
    @{old = "20×64="; new = "90×18="},
    @{old = "24×81="; new = "96×13="},
    @{old = "93×11="; new = "35×68="},
    @{old = "51×83="; new = "82×51="},
    @{old = "92×66="; new = "75×26="},
    @{old = "60×13="; new = "70×93="},
    @{old = "67×37="; new = "21×20="},
    @{old = "97×87="; new = "51×56="},
    @{old = "99×79="; new = "95×16="},
    @{old = "90×34="; new = "32×20="},
    @{old = "44×45="; new = "27×84="},
    @{old = "11×92="; new = "25×66="},
    @{old = "27×78="; new = "84×42="},
    @{old = "79×83="; new = "52×27="},
    @{old = "95×79="; new = "30×33="},
    @{old = "18×22="; new = "24×34="},
    @{old = "36×75="; new = "31×88="},
    @{old = "83×70="; new = "51×52="},
    @{old = "41×26="; new = "82×83="},
    @{old = "21×77="; new = "31×68="},
    @{old = "58×79="; new = "16×55="},
    @{old = "35×89="; new = "99×73="},
    @{old = "49×65="; new = "87×49="},
    @{old = "18×56="; new = "33×41="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
